# Auto-generated edit script applying numeric updates to Anima_Profits workbook
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 10000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10586

$ws.Range("H40").Value = 989
$ws.Range("I40").Value = 1001
$ws.Range("J40").Value = 985.5714
$ws.Range("K40").Value = 1001
$ws.Range("L40").Value = 985.5714
$ws.Range("M40").Value = -826
$ws.Range("N40").Value = -1335.5714

$ws.Range("H92").Value = 25641776
$ws.Range("I92").Value = 27778258
$ws.Range("J92").Value = 4000
$ws.Range("K92").Value = 27778258
$ws.Range("L92").Value = 4000
$ws.Range("M92").Value = -27777010
$ws.Range("N92").Value = -6496

$ws.Range("H106").Value = 20008914
$ws.Range("I106").Value = 25010142
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 25010142
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -25009511
$ws.Range("N106").Value = -5262

$ws.Range("H107").Value = 1028.32
$ws.Range("I107").Value = 1148
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 1148
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 772
$ws.Range("N107").Value = -4240

$ws.Range("H111").Value = 1382.5
$ws.Range("I111").Value = 1172.5714
$ws.Range("K111").Value = 3517.7142
$ws.Range("M111").Value = -450.7142000000003

$ws.Range("H112").Value = 4061.1865
$ws.Range("J112").Value = 4113.9653
$ws.Range("L112").Value = 12341.8959
$ws.Range("N112").Value = -14557.8959

$ws.Range("H118").Value = 930.6923
$ws.Range("I118").Value = 290
$ws.Range("J118").Value = 3066.3333
$ws.Range("K118").Value = 870
$ws.Range("L118").Value = 9198.999899999999
$ws.Range("M118").Value = 787
$ws.Range("N118").Value = -12512.9999

$ws.Range("H125").Value = 2368.5715
$ws.Range("J125").Value = 2368.5715
$ws.Range("L125").Value = 21317.1435
$ws.Range("N125").Value = -26237.1435

$ws.Range("H132").Value = 3663.75
$ws.Range("I132").Value = 3485.2693
$ws.Range("J132").Value = 5984
$ws.Range("K132").Value = 10455.8079
$ws.Range("L132").Value = 17952
$ws.Range("M132").Value = -7925.8079
$ws.Range("N132").Value = -23012

$ws.Range("H137").Value = 1995.3704
$ws.Range("I137").Value = 2187.7273
$ws.Range("K137").Value = 6563.1819
$ws.Range("M137").Value = -4013.1819

$ws.Range("H138").Value = 1892.2959
$ws.Range("I138").Value = 748.88464
$ws.Range("J138").Value = 2305.1943
$ws.Range("K138").Value = 2246.65392
$ws.Range("L138").Value = 6915.5829
$ws.Range("M138").Value = 2893.34608
$ws.Range("N138").Value = -17195.5829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2069.0557
$ws.Range("I45").Value = 1876.8462
$ws.Range("K45").Value = 1876.8462
$ws.Range("M45").Value = -1499.8462

$ws.Range("H74").Value = 19233144
$ws.Range("I74").Value = 1742.4
$ws.Range("J74").Value = 23812050
$ws.Range("K74").Value = 1742.4
$ws.Range("L74").Value = 23812050
$ws.Range("M74").Value = -868.4000000000001
$ws.Range("N74").Value = -23813798

$ws.Range("H77").Value = 19233144
$ws.Range("I77").Value = 1742.4
$ws.Range("J77").Value = 23812050
$ws.Range("K77").Value = 8712
$ws.Range("L77").Value = 119060250
$ws.Range("M77").Value = -4344
$ws.Range("N77").Value = -119068986

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6613.907
$ws.Range("I31").Value = 2055.0625
$ws.Range("K31").Value = 2055.0625
$ws.Range("M31").Value = -1760.0625

$ws.Range("H34").Value = 6613.907
$ws.Range("I34").Value = 2055.0625
$ws.Range("K34").Value = 2055.0625
$ws.Range("M34").Value = -1853.0625

$ws.Range("H132").Value = 2582.7856
$ws.Range("I132").Value = 2158.8572
$ws.Range("J132").Value = 3006.7144
$ws.Range("K132").Value = 6476.571599999999
$ws.Range("L132").Value = 9020.143199999999
$ws.Range("M132").Value = -3946.571599999999
$ws.Range("N132").Value = -14080.1432

$ws.Range("I134").Value = 22737894
$ws.Range("J134").Value = 2158.3333
$ws.Range("K134").Value = 68213682
$ws.Range("L134").Value = 6474.999899999999
$ws.Range("M134").Value = -68211147
$ws.Range("N134").Value = -11544.9999

$ws.Range("H135").Value = 54500
$ws.Range("J135").Value = 54500
$ws.Range("L135").Value = 54500
$ws.Range("N135").Value = -64640

$ws.Range("H140").Value = 78779.664
$ws.Range("J140").Value = 78779.664
$ws.Range("L140").Value = 78779.664
$ws.Range("N140").Value = -89139.664

$ws.Range("H141").Value = 500000.5
$ws.Range("I141").Value = 66667.336
$ws.Range("J141").Value = 1800000
$ws.Range("K141").Value = 66667.336
$ws.Range("L141").Value = 1800000
$ws.Range("M141").Value = -61487.336
$ws.Range("N141").Value = -1810360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 76923270
$ws.Range("I14").Value = 76923270
$ws.Range("K14").Value = 230769810
$ws.Range("M14").Value = -230769637

$ws.Range("H19").Value = 800.3333
$ws.Range("I19").Value = 800.3333
$ws.Range("K19").Value = 2400.9999
$ws.Range("M19").Value = -2226.9999

$ws.Range("H74").Value = 2875
$ws.Range("I74").Value = 2000
$ws.Range("K74").Value = 6000
$ws.Range("M74").Value = -4939

$ws.Range("H77").Value = 2875
$ws.Range("I77").Value = 2000
$ws.Range("K77").Value = 18000
$ws.Range("M77").Value = -12696

$ws.Range("H131").Value = 2800.9155
$ws.Range("I131").Value = 815
$ws.Range("J131").Value = 2858.4783
$ws.Range("K131").Value = 2445
$ws.Range("L131").Value = 8575.4349
$ws.Range("M131").Value = 2595
$ws.Range("N131").Value = -18655.4349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 50780
$ws.Range("J135").Value = 50780
$ws.Range("L135").Value = 50780
$ws.Range("N135").Value = -60920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 38464240
$ws.Range("I82").Value = 45457100
$ws.Range("K82").Value = 45457100
$ws.Range("M82").Value = -45456739

$ws.Range("H85").Value = 38464240
$ws.Range("I85").Value = 45457100
$ws.Range("K85").Value = 45457100
$ws.Range("M85").Value = -45455852

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 56676.332
$ws.Range("I30").Value = 30009
$ws.Range("K30").Value = 30009
$ws.Range("M30").Value = -29902

